$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.093.34"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.955.39"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.99"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4895"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2978"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06862"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.16"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "108.14"
$ws.Range("E11").Value = "  -3.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07757"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.921.02"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.455"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7088"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.41"
$ws.Range("E16").Value = "  -4.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.940.19"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007769"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.181.66"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.497"
$ws.Range("E22").Value = "  -3.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.521"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.19"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.09"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.221"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1055"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.424"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.586"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.589"
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.455"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04964"
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7546"
$ws.Range("E35").Value = "  -3.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.181"
$ws.Range("E36").Value = "  +1.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.732"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02036"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.707"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.176"
$ws.Range("E40").Value = "  +6.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.437"
$ws.Range("E41").Value = "  +8.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4511"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.299"
$ws.Range("E43").Value = "  +12.50%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.32"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8819"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "72.58"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.464"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "965.34"
$ws.Range("E49").Value = "  +6.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1265"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2591"
$ws.Range("E51").Value = "  +1.67%  "
